# Update "想去人数" (F column) values on the "展览" and "全部类型" worksheets
# to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1902
$ws1.Range("F3").Value  = 511
$ws1.Range("F4").Value  = 171
$ws1.Range("F5").Value  = 180
$ws1.Range("F6").Value  = 2726
$ws1.Range("F7").Value  = 190
$ws1.Range("F8").Value  = 97
$ws1.Range("F10").Value = 1576
$ws1.Range("F11").Value = 553
$ws1.Range("F15").Value = 25
$ws1.Range("F17").Value = 8
$ws1.Range("F18").Value = 217
$ws1.Range("F22").Value = 10
$ws1.Range("F23").Value = 217
$ws1.Range("F25").Value = 1744
$ws1.Range("F27").Value = 419
$ws1.Range("F28").Value = 80
$ws1.Range("F30").Value = 216
$ws1.Range("F31").Value = 310
$ws1.Range("F32").Value = 449

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1902
$ws4.Range("F4").Value  = 511
$ws4.Range("F5").Value  = 171
$ws4.Range("F6").Value  = 180
$ws4.Range("F7").Value  = 2726
$ws4.Range("F8").Value  = 190
$ws4.Range("F9").Value  = 97
$ws4.Range("F11").Value = 1576
$ws4.Range("F12").Value = 553
$ws4.Range("F16").Value = 25
$ws4.Range("F18").Value = 8
$ws4.Range("F19").Value = 217
$ws4.Range("F23").Value = 10
$ws4.Range("F24").Value = 217
$ws4.Range("F26").Value = 1744
$ws4.Range("F28").Value = 419
$ws4.Range("F29").Value = 80
$ws4.Range("F31").Value = 216
$ws4.Range("F32").Value = 310
$ws4.Range("F33").Value = 449
